$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.523.11"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.561.92"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.31"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.486"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.76"
$ws.Range("E8").Value = "  +4.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.243"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.787.72"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.560.86"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.565.40"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.512"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.12"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.36"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.35"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0673"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.90"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.95"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.17"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.76"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.103"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.20"
$ws.Range("E29").Value = "  -2.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0460"
$ws.Range("E30").Value = "  -4.55%  "
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.16"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.386.70"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.97"
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.69"
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("E38").Value = "  -2.63%  "
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.516"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.770"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0461"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.70"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.699.77"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("E48").Value = "  -5.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.04"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.84"
$ws.Range("E50").Value = "  +5.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0100"
$ws.Range("E51").Value = "  +0.69%  "
